$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 898.9286
$ws.Range("J129").Value = 911.9074000000001
$ws.Range("L129").Value = 2735.7222
$ws.Range("N129").Value = -12735.7222

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 25451.666
$ws.Range("J24").Value = 25451.666
$ws.Range("L24").Value = 25451.666
$ws.Range("N24").Value = -26199.666
$ws.Range("H28").Value = 8035.5
$ws.Range("I28").Value = 5642.6
$ws.Range("K28").Value = 5642.6
$ws.Range("M28").Value = -5450.6
$ws.Range("H80").Value = 41836.332
$ws.Range("J80").Value = 41836.332
$ws.Range("L80").Value = 41836.332
$ws.Range("N80").Value = -43832.332
$ws.Range("H83").Value = 41836.332
$ws.Range("J83").Value = 41836.332
$ws.Range("L83").Value = 125508.996
$ws.Range("N83").Value = -135492.996
$ws.Range("H92").Value = 24030
$ws.Range("J92").Value = 24030
$ws.Range("L92").Value = 24030
$ws.Range("N92").Value = -29022
$ws.Range("H93").Value = 28433.334
$ws.Range("J93").Value = 28433.334
$ws.Range("L93").Value = 28433.334
$ws.Range("N93").Value = -33425.334
$ws.Range("H94").Value = 25125
$ws.Range("J94").Value = 25125
$ws.Range("L94").Value = 25125
$ws.Range("N94").Value = -26927
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492
$ws.Range("H97").Value = 1463.5
$ws.Range("I97").Value = 1353.1428
$ws.Range("J97").Value = 1849.75
$ws.Range("K97").Value = 1353.1428
$ws.Range("L97").Value = 1849.75
$ws.Range("M97").Value = -857.1428000000001
$ws.Range("N97").Value = -2841.75
$ws.Range("H98").Value = 32556.715
$ws.Range("J98").Value = 32556.715
$ws.Range("L98").Value = 32556.715
$ws.Range("N98").Value = -38546.715
$ws.Range("H99").Value = 8035.5
$ws.Range("I99").Value = 5642.6
$ws.Range("K99").Value = 5642.6
$ws.Range("M99").Value = -2647.6
$ws.Range("H100").Value = 25451.666
$ws.Range("J100").Value = 25451.666
$ws.Range("L100").Value = 25451.666
$ws.Range("N100").Value = -27615.666
$ws.Range("H102").Value = 4516.154
$ws.Range("I102").Value = 3238.75
$ws.Range("K102").Value = 3238.75
$ws.Range("M102").Value = -1616.75
$ws.Range("H104").Value = 38750
$ws.Range("J104").Value = 38750
$ws.Range("L104").Value = 38750
$ws.Range("N104").Value = -45738
$ws.Range("H105").Value = 38246.668
$ws.Range("J105").Value = 38246.668
$ws.Range("L105").Value = 38246.668
$ws.Range("N105").Value = -45234.668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 23333.334
$ws.Range("J92").Value = 23333.334
$ws.Range("L92").Value = 23333.334
$ws.Range("N92").Value = -28325.334
$ws.Range("H93").Value = 37500
$ws.Range("J93").Value = 37500
$ws.Range("L93").Value = 37500
$ws.Range("N93").Value = -41244
$ws.Range("H95").Value = 23017.572
$ws.Range("J95").Value = 23017.572
$ws.Range("L95").Value = 23017.572
$ws.Range("N95").Value = -28509.572
$ws.Range("H96").Value = 35500
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("H97").Value = 9212.5
$ws.Range("I97").Value = 4618.6665
$ws.Range("J97").Value = 16103.25
$ws.Range("K97").Value = 4618.6665
$ws.Range("L97").Value = 16103.25
$ws.Range("M97").Value = -3627.6665
$ws.Range("N97").Value = -18085.25
$ws.Range("H99").Value = 879.6429000000001
$ws.Range("I99").Value = 869.5454999999999
$ws.Range("K99").Value = 869.5454999999999
$ws.Range("M99").Value = 628.4545000000001
$ws.Range("H100").Value = 29682.143
$ws.Range("J100").Value = 29682.143
$ws.Range("L100").Value = 29682.143
$ws.Range("N100").Value = -31846.143
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H102").Value = 9742.857
$ws.Range("I102").Value = 6366.6665
$ws.Range("K102").Value = 6366.6665
$ws.Range("M102").Value = -3121.6665
$ws.Range("H103").Value = 5657
$ws.Range("J103").Value = 5657
$ws.Range("L103").Value = 5657
$ws.Range("N103").Value = -8001
$ws.Range("H105").Value = 2175490.8
$ws.Range("I105").Value = 1463.2142
$ws.Range("J105").Value = 5557311
$ws.Range("K105").Value = 1463.2142
$ws.Range("L105").Value = 5557311
$ws.Range("M105").Value = 283.7858000000001
$ws.Range("N105").Value = -5560805

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1202.25
$ws.Range("I16").Value = 1110.75
$ws.Range("J16").Value = 1385.25
$ws.Range("K16").Value = 1110.75
$ws.Range("L16").Value = 1385.25
$ws.Range("M16").Value = -823.75
$ws.Range("N16").Value = -1959.25
$ws.Range("H43").Value = 31144.5
$ws.Range("J43").Value = 31144.5
$ws.Range("L43").Value = 31144.5
$ws.Range("N43").Value = -31512.5
$ws.Range("H92").Value = 22192
$ws.Range("J92").Value = 22192
$ws.Range("L92").Value = 22192
$ws.Range("N92").Value = -27184
$ws.Range("H93").Value = 18840
$ws.Range("I93").Value = 18600
$ws.Range("J93").Value = 19800
$ws.Range("K93").Value = 18600
$ws.Range("M93").Value = -16728
$ws.Range("N93").Value = -23544
$ws.Range("H95").Value = 15830
$ws.Range("J95").Value = 15830
$ws.Range("L95").Value = 15830
$ws.Range("N95").Value = -21322
$ws.Range("H96").Value = 12722.8
$ws.Range("J96").Value = 12722.8
$ws.Range("L96").Value = 12722.8
$ws.Range("N96").Value = -18214.8
$ws.Range("H101").Value = 31144.5
$ws.Range("J101").Value = 31144.5
$ws.Range("L101").Value = 31144.5
$ws.Range("N101").Value = -37634.5
$ws.Range("H103").Value = 23982.625
$ws.Range("I103").Value = 19500
$ws.Range("J103").Value = 25476.834
$ws.Range("K103").Value = 19500
$ws.Range("L103").Value = 25476.834
$ws.Range("M103").Value = -18328
$ws.Range("N103").Value = -27820.834
$ws.Range("H104").Value = 30650
$ws.Range("J104").Value = 30650
$ws.Range("N104").Value = -35892
$ws.Range("H105").Value = 13890421
$ws.Range("I105").Value = 17857970
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 17857970
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -17856223
$ws.Range("N105").Value = -7494
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H113").Value = 1202.25
$ws.Range("I113").Value = 1110.75
$ws.Range("J113").Value = 1385.25
$ws.Range("K113").Value = 1110.75
$ws.Range("L113").Value = 1385.25
$ws.Range("M113").Value = 1059.25
$ws.Range("N113").Value = -5725.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 929.2857
$ws.Range("J68").Value = 1091.7188
$ws.Range("L68").Value = 3275.1564
$ws.Range("N68").Value = -4897.1564
$ws.Range("H71").Value = 929.2857
$ws.Range("J71").Value = 1091.7188
$ws.Range("L71").Value = 9825.469200000001
$ws.Range("N71").Value = -17937.4692
$ws.Range("H76").Value = 4511.25
$ws.Range("J76").Value = 5013.5
$ws.Range("L76").Value = 15040.5
$ws.Range("N76").Value = -15806.5
$ws.Range("H79").Value = 4511.25
$ws.Range("J79").Value = 5013.5
$ws.Range("L79").Value = 15040.5
$ws.Range("N79").Value = -17692.5
$ws.Range("H131").Value = 800.4
$ws.Range("J131").Value = 826.4316
$ws.Range("L131").Value = 2479.2948
$ws.Range("N131").Value = -12559.2948

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 15750
$ws.Range("J94").Value = 15750
$ws.Range("N94").Value = -17102
$ws.Range("H97").Value = 2139.7222
$ws.Range("I97").Value = 1015.6923
$ws.Range("J97").Value = 5062.2
$ws.Range("K97").Value = 1015.6923
$ws.Range("L97").Value = 5062.2
$ws.Range("M97").Value = -519.6923
$ws.Range("N97").Value = -6054.2
$ws.Range("H100").Value = 40000
$ws.Range("J100").Value = 40000
$ws.Range("L100").Value = 40000
$ws.Range("N100").Value = -42164
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H132").Value = 71744.125
$ws.Range("I132").Value = 17851.666
$ws.Range("J132").Value = 104079.6
$ws.Range("K132").Value = 53554.99800000001
$ws.Range("L132").Value = 312238.8
$ws.Range("M132").Value = -51024.99800000001
$ws.Range("N132").Value = -317298.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1191.1538
$ws.Range("I100").Value = 676
$ws.Range("K100").Value = 1352
$ws.Range("M100").Value = -811
